$d = $word.ActiveDocument

# The document previously split each "<id>...</id>" paragraph into three runs:
#   <id>  (tagged run)  +  p074r_N (plain text run)  +  </id>  (tagged run)
# The edit merges each of these triples into a single run containing the full
# "<id>p074r_N</id>" text (keeping the tag-run formatting: Courier New / 7f6000
# / sz 18 / szCs 18), for ids p074r_1, p074r_2 and p074r_3.

$ids = @("p074r_1", "p074r_2", "p074r_3")

foreach ($id in $ids) {
    $search = "<id>" + $id + "</id>"
    $rng = $d.Content
    $found = $rng.Find.Execute($search, $false, $false, $false, $false, $false, $true, 1, $false, $search, 2)
}
